# Updates cryptos list values (price & volume columns) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: RowNumber, Coin(B), Link(C), Price(D), Volume1h(E)
$data = @(
    ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "30.193.04", "  -0.47%  ")
    ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.858.40", "  -0.60%  ")
    ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "'0.9991", "  -0.13%  ")
    ,@(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "'235.62", "  -0.07%  ")
    ,@(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "'0.9989", "  -0.15%  ")
    ,@(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "'0.4695", "  +0.30%  ")
    ,@(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "'0.2896", "  +1.65%  ")
    ,@(9, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "'0.06596", "  +0.71%  ")
    ,@(10, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "'21.92", "  +1.72%  ")
    ,@(11, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "'0.07983", "  +1.28%  ")
    ,@(12, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "'97.55", "  -0.61%  ")
    ,@(13, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.857.38", "  -0.76%  ")
    ,@(14, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "'5.115", "  +0.19%  ")
    ,@(15, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "'0.6779", "  +0.21%  ")
    ,@(16, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "'268.75", "  -3.17%  ")
    ,@(17, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "30.170.28", "  -0.52%  ")
    ,@(18, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "'13.63", "  +7.10%  ")
    ,@(19, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "'0.000007699", "  +5.27%  ")
    ,@(20, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "'0.9990", "  -0.15%  ")
    ,@(21, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.100.06", "  -0.88%  ")
    ,@(22, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "'0.9988", "  -0.16%  ")
    ,@(23, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "'5.222", "  -4.52%  ")
    ,@(24, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "'6.170", "  +0.25%  ")
    ,@(25, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "'167.10", "  +0.98%  ")
    ,@(26, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "'9.193", "  +0.48%  ")
    ,@(27, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "'18.96", "  -0.91%  ")
    ,@(28, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "'1.943", "  +0.42%  ")
    ,@(29, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "'1.376", "  -0.34%  ")
    ,@(30, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "'0.09898", "  +2.68%  ")
    ,@(31, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "'4.327", "  -1.50%  ")
    ,@(32, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "'1.466", "  -0.59%  ")
    ,@(33, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "'4.036", "  -1.51%  ")
    ,@(34, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "'0.04711", "  +0.19%  ")
    ,@(35, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "'1.123", "  -0.41%  ")
    ,@(36, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "'0.7018", "  -0.64%  ")
    ,@(37, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "'2.703", "  -0.93%  ")
    ,@(38, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "'0.01875", "  +0.85%  ")
    ,@(39, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "'2.607", "  +2.70%  ")
    ,@(40, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "'6.331", "  -0.15%  ")
    ,@(41, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "'73.51", "  -0.90%  ")
    ,@(42, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "'1.936", "  -1.03%  ")
    ,@(43, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "'0.8388", "  -1.32%  ")
    ,@(44, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "'0.9979", "  -0.24%  ")
    ,@(45, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "'0.4152", "  -0.92%  ")
    ,@(46, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "'103.58", "  -0.35%  ")
    ,@(47, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "'7.079", "  -1.76%  ")
    ,@(48, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "'939.15", "  +0.13%  ")
    ,@(49, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "'9.158", "  -1.43%  ")
    ,@(50, "Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "'34.02", "  -0.53%  ")
    ,@(51, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "'0.05655", "  +0.40%  ")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

